# Apply the changes described by the diff to the presentation.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" date field text from
#    12/5/2018 -> 4/2/2019 everywhere it appears: the slide master, all
#    slide layouts and the notes master.
# ---------------------------------------------------------------------
$newDate = "4/2/2019"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "12/5/2018") {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholders $notesMaster.Shapes

# ---------------------------------------------------------------------
# 2. Rename the class-diagram shapes on slide 1 and resize three of the
#    attribute boxes, then drop the now-unused "Address" attribute box
#    and its connector.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes

foreach ($shape in $shapes) {
    if (-not $shape.HasTextFrame) { continue }
    switch ($shape.TextFrame.TextRange.Text) {
        "VersionedAddressBook" { $shape.TextFrame.TextRange.Text = "EntryBook" }
        "UniquePersonList"     { $shape.TextFrame.TextRange.Text = "UniqueEntryList" }
        "Person"                { $shape.TextFrame.TextRange.Text = "Entry" }
        "Name" {
            $shape.TextFrame.TextRange.Text = "Title"
            $shape.Width = 867270 / 12700
        }
        "Phone" {
            $shape.TextFrame.TextRange.Text = "Description"
            $shape.Width = 867270 / 12700
        }
        "Email" {
            $shape.TextFrame.TextRange.Text = "Link"
            $shape.Width = 867270 / 12700
        }
    }
}

# Remove the "Address" attribute box and its connector (Elbow Connector 85).
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shape = $shapes.Item($i)
    $isAddressBox = $shape.HasTextFrame -and ($shape.TextFrame.TextRange.Text -eq "Address")
    $isAddressConnector = ($shape.Name -eq "Elbow Connector 85")
    if ($isAddressBox -or $isAddressConnector) {
        $shape.Delete()
    }
}
